$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (prices/volumes refreshed, two coin rows swapped
# in rank order). Price (column D) and Volume(1h) (column E) are plain
# text cells in this sheet, not numbers - a leading "'" is used on the
# purely-numeric-looking price strings so Excel keeps storing them as
# text (matching the original formatting) instead of silently coercing
# them to a Double and dropping significant trailing zeros.

$ws.Range("D2").Value = "65.530.84"
$ws.Range("E2").Value = "  +0.57%  "

$ws.Range("D3").Value = "3.558.89"
$ws.Range("E3").Value = "  +0.42%  "

$ws.Range("E4").Value = "  +0.30%  "

$ws.Range("D5").Value = "'599.59"
$ws.Range("E5").Value = "  +0.08%  "

$ws.Range("D6").Value = "'136.09"
$ws.Range("E6").Value = "  -1.83%  "

$ws.Range("D7").Value = "3.557.00"
$ws.Range("E7").Value = "  +0.38%  "

$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("D9").Value = "'0.497"
$ws.Range("E9").Value = "  +0.42%  "

$ws.Range("D10").Value = "'0.124"
$ws.Range("E10").Value = "  -0.49%  "

$ws.Range("D11").Value = "'7.14"
$ws.Range("E11").Value = "  +3.27%  "

$ws.Range("D12").Value = "'0.389"
$ws.Range("E12").Value = "  +0.33%  "

$ws.Range("D13").Value = "4.191.28"
$ws.Range("E13").Value = "  +1.08%  "

$ws.Range("D14").Value = "'27.71"
$ws.Range("E14").Value = "  +1.12%  "

$ws.Range("D15").Value = "'0.0000183"
$ws.Range("E15").Value = "  -0.65%  "

$ws.Range("D16").Value = "3.593.27"
$ws.Range("E16").Value = "  +1.42%  "

$ws.Range("E17").Value = "  -0.01%  "

$ws.Range("D18").Value = "64.788.19"
$ws.Range("E18").Value = "  -0.48%  "

$ws.Range("D19").Value = "'9.82"
$ws.Range("E19").Value = "  -2.57%  "

$ws.Range("D20").Value = "'14.54"
$ws.Range("E20").Value = "  +2.08%  "

$ws.Range("D21").Value = "'5.77"
$ws.Range("E21").Value = "  -2.15%  "

$ws.Range("D22").Value = "'393.76"
$ws.Range("E22").Value = "  +0.38%  "

$ws.Range("D23").Value = "'0.582"
$ws.Range("E23").Value = "  +1.03%  "

$ws.Range("D24").Value = "3.711.68"
$ws.Range("E24").Value = "  +0.74%  "

$ws.Range("E25").Value = "  +1.53%  "

$ws.Range("D26").Value = "'0.996"
$ws.Range("E26").Value = "  -0.34%  "

$ws.Range("D27").Value = "'0.0000116"
$ws.Range("E27").Value = "  +1.85%  "

$ws.Range("D28").Value = "'7.92"
$ws.Range("E28").Value = "  +0.12%  "

$ws.Range("E29").Value = "  +17.06%  "

$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.33%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'8.54"
$ws.Range("E31").Value = "  +2.10%  "

$ws.Range("D32").Value = "'2.31"
$ws.Range("E32").Value = "  +1.18%  "

$ws.Range("D33").Value = "3.580.32"
$ws.Range("E33").Value = "  +0.46%  "

$ws.Range("D34").Value = "'24.29"
$ws.Range("E34").Value = "  +1.66%  "

$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.02%  "

$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "'0.147"
$ws.Range("E36").Value = "  +1.21%  "

$ws.Range("D37").Value = "'5.33"
$ws.Range("E37").Value = "  +5.38%  "

$ws.Range("D38").Value = "'1.59"
$ws.Range("E38").Value = "  -0.04%  "

$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").Value = "'169.43"
$ws.Range("E39").Value = "  -0.29%  "

$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").Value = "'6.91"
$ws.Range("E40").Value = "  -0.74%  "

$ws.Range("D41").Value = "'0.0829"
$ws.Range("E41").Value = "  +2.58%  "

$ws.Range("D42").Value = "'0.830"
$ws.Range("E42").Value = "  +0.70%  "

$ws.Range("D43").Value = "'26.18"
$ws.Range("E43").Value = "  -2.23%  "

$ws.Range("E44").Value = "  +4.25%  "

$ws.Range("D45").Value = "'42.83"
$ws.Range("E45").Value = "  +0.72%  "

$ws.Range("E46").Value = "  +0.44%  "

$ws.Range("D47").Value = "'4.47"
$ws.Range("E47").Value = "  +0.25%  "

$ws.Range("D48").Value = "'1.68"
$ws.Range("E48").Value = "  +0.40%  "

$ws.Range("D49").Value = "'6.89"
$ws.Range("E49").Value = "  +0.56%  "

$ws.Range("D50").Value = "2.430.93"
$ws.Range("E50").Value = "  +0.82%  "

$ws.Range("D51").Value = "'0.900"
$ws.Range("E51").Value = "  +5.47%  "
